# Update workbook to match scraped data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Ticket-interest counter bumps that happened before the new row is inserted.
$ws1.Cells.Item(3,6).Value = 2464
$ws1.Cells.Item(5,6).Value = 1687
$ws1.Cells.Item(7,6).Value = 317
$ws1.Cells.Item(8,6).Value = 615
$ws1.Cells.Item(9,6).Value = 3520
$ws1.Cells.Item(10,6).Value = 920
$ws1.Cells.Item(12,6).Value = 1567
$ws1.Cells.Item(14,6).Value = 886
$ws1.Cells.Item(15,6).Value = 14
$ws1.Cells.Item(16,6).Value = 1279
$ws1.Cells.Item(17,6).Value = 1785
$ws1.Cells.Item(19,6).Value = 449
$ws1.Cells.Item(20,6).Value = 1545

# A brand-new event was scraped in between the existing row 20 and row 21,
# so insert a fresh row at position 21 (pushes old rows 21-28 down to 22-29).
$ws1.Rows.Item(21).Insert()

# Insert() does not restyle the freshly-created row; copy the index column's
# formatting down from the row above so A21 keeps the bold/bordered look.
$ws1.Cells.Item(20,1).Copy()
$ws1.Cells.Item(21,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Cells.Item(21,1).Value = 20
$ws1.Cells.Item(21,2).Value = "'2024-11-03"
$ws1.Cells.Item(21,3).Value = "上海·恋与深空契约3.0同人Only回溯之时"
$ws1.Cells.Item(21,4).Value = "南翔镇嘉前路275号 嘉美莉雅一站式宴会中心"
$ws1.Cells.Item(21,5).Value = "2024.11.03 10:30-11.03 21:00"
$ws1.Cells.Item(21,6).Value = 2
$ws1.Cells.Item(21,7).Value = 98
$ws1.Cells.Item(21,8).Value = "https://show.bilibili.com/platform/detail.html?id=93449"
$ws1.Cells.Item(21,9).Value = "//i2.hdslb.com/bfs/openplatform/202410/e6KxRlyj1728566600416.jpeg"

# Ticket-interest counter bumps on the rows that shifted down one slot.
$ws1.Cells.Item(23,6).Value = 2089   # 创造力动漫游戏嘉年华2.0
$ws1.Cells.Item(24,6).Value = 4      # 夜蓝诗2.0·恋与深空同人only
$ws1.Cells.Item(26,6).Value = 4260   # 趣元界·第三届ICG动漫游戏博览会
$ws1.Cells.Item(28,6).Value = 2695   # 第二届iPR动漫游戏嘉年华

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(2,6).Value = 81
$ws2.Cells.Item(35,6).Value = 13
$ws2.Cells.Item(47,6).Value = 27
$ws2.Cells.Item(48,6).Value = 27

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Cells.Item(4,6).Value = 2540
$ws3.Cells.Item(6,6).Value = 9573
$ws3.Cells.Item(11,6).Value = 2960
$ws3.Cells.Item(12,6).Value = 464
$ws3.Cells.Item(13,6).Value = 798
$ws3.Cells.Item(14,6).Value = 203

# A new row is appended at the end of the used range (row 15).
$ws3.Cells.Item(14,1).Copy()
$ws3.Cells.Item(15,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Cells.Item(15,1).Value = 14
$ws3.Cells.Item(15,2).Value = "'2024-11-15"
$ws3.Cells.Item(15,3).Value = "上海·「WIND BREAKER × animate cafe」"
$ws3.Cells.Item(15,4).Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws3.Cells.Item(15,5).Value = "2024.11.15 00:00-12.15 23:59"
$ws3.Cells.Item(15,6).Value = 50
$ws3.Cells.Item(15,7).Value = 30
$ws3.Cells.Item(15,8).Value = "https://show.bilibili.com/platform/detail.html?id=93422"
$ws3.Cells.Item(15,9).Value = "//i0.hdslb.com/bfs/openplatform/202410/TGPx1EZW1728892799830.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" (all categories) - an independently generated, static
# snapshot that mirrors counters from the sheets above (it is not resized).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2,6).Value = 2540
$ws4.Cells.Item(6,6).Value = 2464
$ws4.Cells.Item(9,6).Value = 2960
$ws4.Cells.Item(10,6).Value = 798
$ws4.Cells.Item(11,6).Value = 203
$ws4.Cells.Item(12,6).Value = 81
$ws4.Cells.Item(14,6).Value = 1687
$ws4.Cells.Item(15,6).Value = 317
$ws4.Cells.Item(16,6).Value = 615
$ws4.Cells.Item(17,6).Value = 920
$ws4.Cells.Item(19,7).Value = "不可售"
$ws4.Cells.Item(20,6).Value = 886
$ws4.Cells.Item(31,6).Value = 1785
$ws4.Cells.Item(33,6).Value = 1545
$ws4.Cells.Item(40,6).Value = 2089
$ws4.Cells.Item(44,6).Value = 13
$ws4.Cells.Item(49,6).Value = 27
$ws4.Cells.Item(50,6).Value = 27
